# Applies the "Updated symbol list" data refresh (Wed Dec 28 05:08:54 UTC 2022,
# GitHub Actions) to the crypto price table on Sheet1: refreshed prices, the
# "Hora" (hour) column bumped from 4 to 5, a few "Best/Worst in 24h" labels
# moved to a different row, and several coin rows swapped position after the
# upstream re-sort (name/link/price/volume-label move together; row index in
# column A and the scrape date in column F stay put).
#
# Price/Hora cells are stored as genuine text (scraped strings), not numbers -
# e.g. "244.90" must stay "244.90", not become the number 244.9. Assigning a
# numeric-looking string straight to .Value lets Excel auto-convert it to a
# real number, so those values are written with a leading apostrophe, which
# forces text entry exactly like typing into the cell in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''244.90'
$ws.Range("G2").Value = '''5'
# Row 3
$ws.Range("D3").Value = '''23.73'
$ws.Range("G3").Value = '''5'
# Row 4
$ws.Range("D4").Value = '''5.331'
$ws.Range("G4").Value = '''5'
# Row 5
$ws.Range("G5").Value = '''5'
# Row 6
$ws.Range("D6").Value = '''6.479'
$ws.Range("G6").Value = '''5'
# Row 7
$ws.Range("D7").Value = '''3.334'
$ws.Range("G7").Value = '''5'
# Row 8
$ws.Range("D8").Value = '''0.8091'
$ws.Range("G8").Value = '''5'
# Row 9
$ws.Range("D9").Value = '''0.8859'
$ws.Range("G9").Value = '''5'
# Row 10
$ws.Range("D10").Value = '''0.1393'
$ws.Range("G10").Value = '''5'
# Row 11
$ws.Range("D11").Value = '''0.07359'
$ws.Range("G11").Value = '''5'
# Row 12
$ws.Range("D12").Value = '''0.03080'
$ws.Range("G12").Value = '''5'
# Row 13
$ws.Range("D13").Value = '''0.03063'
$ws.Range("G13").Value = '''5'
# Row 14
$ws.Range("D14").Value = '''0.09341'
$ws.Range("G14").Value = '''5'
# Row 15
$ws.Range("D15").Value = '''3.852'
$ws.Range("G15").Value = '''5'
# Row 16
$ws.Range("D16").Value = '''0.001562'
$ws.Range("G16").Value = '''5'
# Row 17
$ws.Range("D17").Value = '''0.04718'
$ws.Range("G17").Value = '''5'
# Row 18
$ws.Range("D18").Value = '''0.0006015'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("G18").Value = '''5'
# Row 19
$ws.Range("D19").Value = '''0.005989'
$ws.Range("G19").Value = '''5'
# Row 20
$ws.Range("D20").Value = '''0.001295'
$ws.Range("G20").Value = '''5'
# Row 21
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = '''0.004654'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("G21").Value = '''5'
# Row 22
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D22").Value = '''0.00008805'
$ws.Range("E22").Value = '21NitroExNTXBestin24h'
$ws.Range("G22").Value = '''5'
# Row 23
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D23").Value = '''3.583'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("G23").Value = '''5'
# Row 24
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D24").Value = '''2.141'
$ws.Range("E24").Value = '23BTSETokenBTSE'
$ws.Range("G24").Value = '''5'
# Row 25
$ws.Range("B25").Value = 'BitpandaEcosystemToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D25").Value = '''0.3179'
$ws.Range("E25").Value = '24BitpandaEcosystemTokenBEST'
$ws.Range("G25").Value = '''5'
# Row 26
$ws.Range("B26").Value = 'ProBitToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D26").Value = '''0.1318'
$ws.Range("E26").Value = '25ProBitTokenPROB'
$ws.Range("G26").Value = '''5'
# Row 27
$ws.Range("B27").Value = 'AAXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("D27").Value = '''0.2000'
$ws.Range("E27").Value = '26AAXTokenAAB'
$ws.Range("G27").Value = '''5'
# Row 28
$ws.Range("B28").Value = 'UpBots'
$ws.Range("C28").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D28").Value = '''0.0002352'
$ws.Range("E28").Value = '27UpBotsUBXT'
$ws.Range("G28").Value = '''5'
# Row 29
$ws.Range("G29").Value = '''5'
# Row 30
$ws.Range("G30").Value = '''5'
# Row 31
$ws.Range("G31").Value = '''5'
# Row 32
$ws.Range("G32").Value = '''5'
# Row 33
$ws.Range("G33").Value = '''5'
# Row 34
$ws.Range("G34").Value = '''5'
# Row 35
$ws.Range("G35").Value = '''5'
# Row 36
$ws.Range("G36").Value = '''5'
# Row 37
$ws.Range("G37").Value = '''5'
# Row 38
$ws.Range("G38").Value = '''5'
# Row 39
$ws.Range("G39").Value = '''5'
# Row 40
$ws.Range("D40").Value = '''0.03814'
$ws.Range("G40").Value = '''5'
# Row 41
$ws.Range("D41").Value = '''0.006353'
$ws.Range("E41").Value = '40KickTokenKICK'
$ws.Range("G41").Value = '''5'
# Row 42
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '''0.1054'
$ws.Range("E42").Value = '41BKEXTokenBKK'
$ws.Range("G42").Value = '''5'
# Row 43
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '''0.002801'
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("G43").Value = '''5'
# Row 44
$ws.Range("D44").Value = '''0.007605'
$ws.Range("G44").Value = '''5'
# Row 45
$ws.Range("D45").Value = '''0.00005457'
$ws.Range("G45").Value = '''5'
# Row 46
$ws.Range("D46").Value = '''0.00000000750'
$ws.Range("G46").Value = '''5'
# Row 47
$ws.Range("D47").Value = '''0.5505'
$ws.Range("G47").Value = '''5'
# Row 48
$ws.Range("D48").Value = '''0.001839'
$ws.Range("G48").Value = '''5'
# Row 49
$ws.Range("D49").Value = '''0.00002101'
$ws.Range("G49").Value = '''5'
# Row 50
$ws.Range("D50").Value = '''0.0002001'
$ws.Range("G50").Value = '''5'
# Row 51
$ws.Range("G51").Value = '''5'
